$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# 1) RightsStatement sheet: update the "note" text in row 2 (E2)
# -----------------------------------------------------------------
$wsRights = $wb.Worksheets.Item("RightsStatement")
$wsRights.Range("E2").Value = "Unless expressly stated otherwise, the organization that has made this Item available makes no warranties about the Item and cannot guarantee the accuracy of this Rights Statement. You are responsible for your own use."

# -----------------------------------------------------------------
# 2) Property sheet: insert new "cms:propertyFilterable" column
#    (B) and append a new "searchable" column (F)
# -----------------------------------------------------------------
$wsProperty = $wb.Worksheets.Item("Property")

# Shift group/label/range columns one to the right, making room for
# the new "cms:propertyFilterable" column in B.
$wsProperty.Range("B:B").Insert(-4161)

$wsProperty.Range("B1").Value = "cms:propertyFilterable"
$wsProperty.Range("F1").Value = "searchable"

$filterableValue = "<class 'filter'>"
for ($r = 2; $r -le 12; $r++) {
    $wsProperty.Cells.Item($r, 2).Value = $filterableValue
    $wsProperty.Cells.Item($r, 6).Value = "'true"
}

# -----------------------------------------------------------------
# 3) Person sheet: change two "relation" URLs from wikidata to
#    wikipedia (rows 3 and 6 -- person1 and person4)
# -----------------------------------------------------------------
$wsPerson = $wb.Worksheets.Item("Person")
$wsPerson.Range("F3").Value = "http://en.wikipedia.org/wiki/Alan_Turing"
$wsPerson.Range("F6").Value = "http://en.wikipedia.org/wiki/Alan_Turing"
